# Final touches for C23A
#
# The test-layout grid on Sheet1 has several pairs of adjacent cells whose
# "door direction test" marker (the darker/"#"-suffixed shared string +
# distinct fill) needs to be swapped with its neighboring base-room cell.
# Each pair below is a full swap of BOTH the cell value and its style
# (fill/font/border/number format), so we round-trip each pair through a
# scratch cell using Copy + PasteSpecial(xlPasteFormats) for the style and
# a plain Value assignment for the content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$scratch = "ZZ1"

function Swap-CellFull($ref1, $ref2) {
    # scratch = full copy (value + style) of ref1
    $ws.Range($ref1).Copy()
    $ws.Range($scratch).PasteSpecial($xlPasteFormats)
    $ws.Range($scratch).Value = $ws.Range($ref1).Value2

    # ref1 = full copy (value + style) of ref2
    $ws.Range($ref2).Copy()
    $ws.Range($ref1).PasteSpecial($xlPasteFormats)
    $ws.Range($ref1).Value = $ws.Range($ref2).Value2

    # ref2 = full copy (value + style) of scratch (== original ref1)
    $ws.Range($scratch).Copy()
    $ws.Range($ref2).PasteSpecial($xlPasteFormats)
    $ws.Range($ref2).Value = $ws.Range($scratch).Value2

    # tidy up the scratch cell
    $ws.Range($scratch).Clear()
}

Swap-CellFull "D3" "E4"
Swap-CellFull "M4" "L8"
Swap-CellFull "D15" "J17"
Swap-CellFull "AC20" "AD21"
Swap-CellFull "D22" "D24"
Swap-CellFull "O22" "O24"
Swap-CellFull "W22" "X22"

# Update the saved cursor / active-cell selection shown when the sheet is
# reopened.
$ws.Range("Q17").Select()
